$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab label shown in workbook.xml's <sheet name="..."> element
$ws.Name = "AlphaFiberF"

# Match row 15's formatting for the first column (style index 1 = bordered/bold/centered)
# by copying the formatted source cell onto the new cell before writing values.
$ws.Range("A15").Copy($ws.Range("A16"))

# Append a new row (row 16) of averaged intensity data, mirroring row 15's
# structure (same "HexGrid-60degTilt5degRes" label, next index value).
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9906828493086272
$ws.Range("D16").Value = 0.9997385364559017
$ws.Range("E16").Value = 0.991866189525917
$ws.Range("F16").Value = 0.9906828493086272
$ws.Range("G16").Value = 1.001583232288522
$ws.Range("H16").Value = 0.9882714118751865
$ws.Range("I16").Value = 0.9929119283421032
$ws.Range("J16").Value = 0.9997385364559017
$ws.Range("K16").Value = 0.9958023629909094
$ws.Range("L16").Value = 0.9932426061497682
$ws.Range("M16").Value = 0.9941756912993762
